# Sprint.0.ASimpleTODOLIST.TRD.xlsx - add "Diseño" and "Recursos" content,
# rename the three sheets and switch the active tab to "Diseño".

$wb = $excel.ActiveWorkbook

# --- Rename sheets -----------------------------------------------------
$wsTRD = $wb.Worksheets.Item(1)
$wsDiseno = $wb.Worksheets.Item(2)
$wsRecursos = $wb.Worksheets.Item(3)

$wsTRD.Name = "TRD"
$wsDiseno.Name = "Diseño"
$wsRecursos.Name = "Recursos"

# --- "Diseño" sheet content ---------------------------------------------
# Filled in the same order the original author typed it: the three
# Plataforma/Lenguaje/MySQL rows first, then the "Tipo de Aplicacion"
# header row above them, then "Base de Datos:", then the server/hosting
# rows - this keeps the shared-string table in the same append order.
$wsDiseno.Range("A3").Value = "Plataforma:"
$wsDiseno.Range("B3").Value = "Apache"
$wsDiseno.Range("A4").Value = "Lenguaje :"
$wsDiseno.Range("B4").Value = "PHP"
$wsDiseno.Range("B5").Value = "MySQL"
$wsDiseno.Range("A2").Value = "Tipo de Aplicacion"
$wsDiseno.Range("B2").Value = "Web"
$wsDiseno.Range("A5").Value = "Base de Datos:"
$wsDiseno.Range("A6").Value = "Numero de Servidores"
$wsDiseno.Range("B6").Value = 1
$wsDiseno.Range("A7").Value = "Tipo de Hosting"
$wsDiseno.Range("B7").Value = "Compartido"

$wsDiseno.Columns.Item(1).ColumnWidth = 18.333333333333332
$wsDiseno.Columns.Item(2).ColumnWidth = 10.333333333333332

# --- "Recursos" sheet content --------------------------------------------
$wsRecursos.Range("B3").Value = "Incluir aqui la suma de horas a desarrollar"
$wsRecursos.Range("A4").Value = "Costo por hora hombre"
$wsRecursos.Range("B4").Value = "Incluir cuando debe de ganarse por programar"
$wsRecursos.Range("A5").Value = "Licencias y paquetes"
$wsRecursos.Range("B5").Value = "Si compran alguna libreria o paquete"
$wsRecursos.Range("A3").Value = "Costo de Desarrollo en Horas Hombre"

$wsRecursos.Range("A7").Value = "Instalacion de producto"
$wsRecursos.Range("A8").Value = "Costo de Operacion"
$wsRecursos.Range("B7").Value = "Estimar que tanto tiempo tardan en instalar el producto y como. Te mando un correo, te administro el hosting etc."
$wsRecursos.Range("B8").Value = "Cuanto cuesta el shared hosting, incluir precio de dominio, hosting y cualquier otro gasto recurrente."

$wsRecursos.Range("A11").Value = "Total"
$wsRecursos.Range("B11").Value = "Sumar el Primer pago"
$wsRecursos.Range("A12").Value = "Anticipo"
$wsRecursos.Range("A13").Value = "Resto "

$wsRecursos.Columns.Item(1).ColumnWidth = 31.5
$wsRecursos.Columns.Item(2).ColumnWidth = 95.0

# --- Selections / active tab --------------------------------------------
# Leave a selection parked on Recursos (not the active sheet) ...
$wsRecursos.Activate()
$wsRecursos.Range("A8").Select() | Out-Null

# ... then activate Diseño last so it becomes the workbook's active tab,
# with B8 selected (one row below the last used row).
$wsDiseno.Activate()
$wsDiseno.Range("B8").Select() | Out-Null

Write-Output "done"
